$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 2025 season rows (201-210) below the existing 2006-2024 data.
# Columns: A=Yr, B=Person, C=SRank, D=Points, E=Bonus, F=PointsBonus,
#          G=Chips, H=Winnings, I=Takehome, J=PersStatus, K=pers_personid

$rows = @(
    @(2025, "Andy",     1, 12, 0, 12, 44600, 60, 40,  "Active", 349),
    @(2025, "Prashant", 2,  9, 0,  9, 27800, 20,  0,  "Active", 365),
    @(2025, "Matt",     3,  9, 0,  9, 27050, 50, 30,  "Active", 362),
    @(2025, "Richard",  4,  6, 0,  6, 20800, 20,  0,  "Active", 366),
    @(2025, "Pepe",     5,  6, 0,  6, 20250,  0, -20, "Active", 364),
    @(2025, "Maisy",    6,  5, 0,  5, 18350,  0, -20, "Active", 360),
    @(2025, "Mark",     7,  5, 0,  5, 17050, 10,  0,  "Active", 361),
    @(2025, "Jon",      8,  4, 0,  4, 13050,  0, -20, "Active", 357),
    @(2025, "Anthony",  9,  4, 0,  4,  9300, 10,  0,  "Active", 350),
    @(2025, "Alex",    10,  0, 0,  0,  4500,  0, -10, "Active", 348)
)

$startRow = 201
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
    $ws.Cells.Item($r, 10).Value = $data[9]
    $ws.Cells.Item($r, 11).Value = $data[10]
}
